# Update "last_edited_time" column (D) values in the
# "Lũy kế ngày LONG XUYÊN" sheet.
#
# Net effect required by the diff:
#   D2:D15   "2024-08-03T03:28:00.000Z" -> "2024-08-03T03:54:00.000Z"
#   D16:D30  "2024-08-03T03:29:00.000Z" -> "2024-08-03T03:54:00.000Z"
#   D54:D84  "2024-08-03T03:30:00.000Z" -> "2024-08-03T03:55:00.000Z"
#   D85:D94  "2024-08-03T03:30:00.000Z" -> "2024-08-03T03:56:00.000Z"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2:D15").Value  = "2024-08-03T03:54:00.000Z"
$ws.Range("D16:D30").Value = "2024-08-03T03:54:00.000Z"
$ws.Range("D54:D84").Value = "2024-08-03T03:55:00.000Z"
$ws.Range("D85:D94").Value = "2024-08-03T03:56:00.000Z"
